$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source dataset gained six new weekly price rows for "Terminal
# Hortofruticola Agro Chillan - Palta" (fecha 2023-04-25 / serial 45041).
# They are inserted at row 775, pushing the existing rows 775:802 down to
# 781:808 (dimension grows from A1:T802 to A1:T808).
$ws.Rows("775:780").Insert()

$ws.Cells.Item(775, 1).Value = 7
$ws.Cells.Item(775, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(775, 3).Value = "Ñuble"
$ws.Cells.Item(775, 4).Value = 45041
$ws.Cells.Item(775, 5).Value = 16
$ws.Cells.Item(775, 6).Value = "Fruta"
$ws.Cells.Item(775, 7).Value = 100106
$ws.Cells.Item(775, 8).Value = "Oleaginosos"
$ws.Cells.Item(775, 9).Value = 100106002
$ws.Cells.Item(775, 10).Value = "Palta"
$ws.Cells.Item(775, 11).Value = "Hass"
$ws.Cells.Item(775, 12).Value = "Especial"
$ws.Cells.Item(775, 13).Value = 80
$ws.Cells.Item(775, 14).Value = 35000
$ws.Cells.Item(775, 15).Value = 35000
$ws.Cells.Item(775, 16).Value = 35000
$ws.Cells.Item(775, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(775, 18).Value = "Perú"
$ws.Cells.Item(775, 19).Value = 3500
$ws.Cells.Item(775, 20).Value = 10


$ws.Cells.Item(776, 1).Value = 7
$ws.Cells.Item(776, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(776, 3).Value = "Ñuble"
$ws.Cells.Item(776, 4).Value = 45041
$ws.Cells.Item(776, 5).Value = 16
$ws.Cells.Item(776, 6).Value = "Fruta"
$ws.Cells.Item(776, 7).Value = 100106
$ws.Cells.Item(776, 8).Value = "Oleaginosos"
$ws.Cells.Item(776, 9).Value = 100106002
$ws.Cells.Item(776, 10).Value = "Palta"
$ws.Cells.Item(776, 11).Value = "Hass"
$ws.Cells.Item(776, 12).Value = "Especial"
$ws.Cells.Item(776, 13).Value = 100
$ws.Cells.Item(776, 14).Value = 5000
$ws.Cells.Item(776, 15).Value = 5000
$ws.Cells.Item(776, 16).Value = 5000
$ws.Cells.Item(776, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(776, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(776, 19).Value = 5000
$ws.Cells.Item(776, 20).Value = 1


$ws.Cells.Item(777, 1).Value = 7
$ws.Cells.Item(777, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(777, 3).Value = "Ñuble"
$ws.Cells.Item(777, 4).Value = 45041
$ws.Cells.Item(777, 5).Value = 16
$ws.Cells.Item(777, 6).Value = "Fruta"
$ws.Cells.Item(777, 7).Value = 100106
$ws.Cells.Item(777, 8).Value = "Oleaginosos"
$ws.Cells.Item(777, 9).Value = 100106002
$ws.Cells.Item(777, 10).Value = "Palta"
$ws.Cells.Item(777, 11).Value = "Hass"
$ws.Cells.Item(777, 12).Value = "Primera"
$ws.Cells.Item(777, 13).Value = 100
$ws.Cells.Item(777, 14).Value = 30000
$ws.Cells.Item(777, 15).Value = 30000
$ws.Cells.Item(777, 16).Value = 30000
$ws.Cells.Item(777, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(777, 18).Value = "Perú"
$ws.Cells.Item(777, 19).Value = 3000
$ws.Cells.Item(777, 20).Value = 10


$ws.Cells.Item(778, 1).Value = 7
$ws.Cells.Item(778, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(778, 3).Value = "Ñuble"
$ws.Cells.Item(778, 4).Value = 45041
$ws.Cells.Item(778, 5).Value = 16
$ws.Cells.Item(778, 6).Value = "Fruta"
$ws.Cells.Item(778, 7).Value = 100106
$ws.Cells.Item(778, 8).Value = "Oleaginosos"
$ws.Cells.Item(778, 9).Value = 100106002
$ws.Cells.Item(778, 10).Value = "Palta"
$ws.Cells.Item(778, 11).Value = "Hass"
$ws.Cells.Item(778, 12).Value = "Primera"
$ws.Cells.Item(778, 13).Value = 100
$ws.Cells.Item(778, 14).Value = 4500
$ws.Cells.Item(778, 15).Value = 4500
$ws.Cells.Item(778, 16).Value = 4500
$ws.Cells.Item(778, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(778, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(778, 19).Value = 4500
$ws.Cells.Item(778, 20).Value = 1


$ws.Cells.Item(779, 1).Value = 7
$ws.Cells.Item(779, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(779, 3).Value = "Ñuble"
$ws.Cells.Item(779, 4).Value = 45041
$ws.Cells.Item(779, 5).Value = 16
$ws.Cells.Item(779, 6).Value = "Fruta"
$ws.Cells.Item(779, 7).Value = 100106
$ws.Cells.Item(779, 8).Value = "Oleaginosos"
$ws.Cells.Item(779, 9).Value = 100106002
$ws.Cells.Item(779, 10).Value = "Palta"
$ws.Cells.Item(779, 11).Value = "Hass"
$ws.Cells.Item(779, 12).Value = "Segunda"
$ws.Cells.Item(779, 13).Value = 80
$ws.Cells.Item(779, 14).Value = 28000
$ws.Cells.Item(779, 15).Value = 28000
$ws.Cells.Item(779, 16).Value = 28000
$ws.Cells.Item(779, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(779, 18).Value = "Perú"
$ws.Cells.Item(779, 19).Value = 2800
$ws.Cells.Item(779, 20).Value = 10


$ws.Cells.Item(780, 1).Value = 7
$ws.Cells.Item(780, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(780, 3).Value = "Ñuble"
$ws.Cells.Item(780, 4).Value = 45041
$ws.Cells.Item(780, 5).Value = 16
$ws.Cells.Item(780, 6).Value = "Fruta"
$ws.Cells.Item(780, 7).Value = 100106
$ws.Cells.Item(780, 8).Value = "Oleaginosos"
$ws.Cells.Item(780, 9).Value = 100106002
$ws.Cells.Item(780, 10).Value = "Palta"
$ws.Cells.Item(780, 11).Value = "Hass"
$ws.Cells.Item(780, 12).Value = "Segunda"
$ws.Cells.Item(780, 13).Value = 80
$ws.Cells.Item(780, 14).Value = 3800
$ws.Cells.Item(780, 15).Value = 3800
$ws.Cells.Item(780, 16).Value = 3800
$ws.Cells.Item(780, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(780, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(780, 19).Value = 3800
$ws.Cells.Item(780, 20).Value = 1
